$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F5").Value = 5003
$ws1.Range("F13").Value = 1394
$ws1.Range("F19").Value = 2656
$ws1.Range("F28").Value = 268

$ws4.Range("F6").Value = 5003
$ws4.Range("F14").Value = 1394
$ws4.Range("F20").Value = 2656
$ws4.Range("F29").Value = 268
